$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a numeric-looking text value (price / percent) while preserving it
# as literal text (matching the source workbook, where these columns are stored
# as plain strings, not numbers) and without leaving a residual Text number format
# on the cell (ClearFormats drops the style index Excel adds for NumberFormat="@").
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "286.09"
Set-TextValue "E2" "1.49%"
Set-TextValue "D3" "29.30"
Set-TextValue "E3" "3.13%"
Set-TextValue "D4" "5.069"
Set-TextValue "E4" "0.36%"
Set-TextValue "D5" "0.06737"
Set-TextValue "E5" "4.19%"
Set-TextValue "D6" "7.337"
Set-TextValue "E6" "1.62%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D7" "1.381"
Set-TextValue "E7" "-0.60%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9011"
Set-TextValue "E8" "-2.97%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1583"
Set-TextValue "E9" "2.96%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.06886"
Set-TextValue "E10" "7.74%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.07601"
Set-TextValue "E11" "0.73%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.02924"
Set-TextValue "E12" "0.73%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.08987"
Set-TextValue "E13" "0.31%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001579"
Set-TextValue "E14" "-1.05%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D15" "0.04479"
Set-TextValue "E15" "1.61%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D16" "0.0006473"
Set-TextValue "E16" "0.81%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D17" "0.006606"
Set-TextValue "E17" "8.20%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.453"
Set-TextValue "E18" "0.36%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D19" "3.435"
Set-TextValue "E19" "1.55%"
Set-TextValue "D20" "2.230"
Set-TextValue "E20" "0.00%"
Set-TextValue "D21" "0.3206"
Set-TextValue "E21" "0.55%"
Set-TextValue "E22" "2.11%"
Set-TextValue "D23" "4.002"
Set-TextValue "E23" "-1.35%"
Set-TextValue "E24" "1.83%"
Set-TextValue "D25" "0.001203"
Set-TextValue "E25" "1.30%"
Set-TextValue "D26" "0.004374"
Set-TextValue "E26" "8.38%"
Set-TextValue "D27" "0.0001167"
Set-TextValue "E27" "-6.85%"
Set-TextValue "D28" "0.0001615"
Set-TextValue "E28" "-0.72%"
Set-TextValue "D40" "0.04239"
Set-TextValue "E40" "3.20%"
Set-TextValue "D41" "0.006812"
Set-TextValue "E41" "5.93%"
Set-TextValue "D42" "0.1238"
Set-TextValue "E42" "1.83%"
Set-TextValue "D43" "0.002185"
Set-TextValue "E43" "3.79%"
Set-TextValue "D44" "0.01154"
Set-TextValue "E44" "-4.25%"
Set-TextValue "D45" "0.00005736"
Set-TextValue "E45" "1.47%"
Set-TextValue "D46" "1.926"
Set-TextValue "E46" "-2.01%"
Set-TextValue "E47" "14.80%"
